$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.830195666666667
$ws.Range("H2").Value = 14.490587
$ws.Range("I2").Value = 0.1644281803937369
$ws.Range("J2").Value = 0.1663774220310514
$ws.Range("M2").Value = 0.013841
$ws.Range("N2").Value = 0.041523
$ws.Range("O2").Value = 0.001379248562145083
$ws.Range("P2").Value = 0.001433331793786624
$ws.Range("Q2").Value = 0.06685473822233333
$ws.Range("R2").Value = 0.601692644001
$ws.Range("S2").Value = 0.0002267873313841939
$ws.Range("T2").Value = 0.0002384740487653611
$ws.Range("G3").Value = 4.830195666666667
$ws.Range("H3").Value = 14.490587
$ws.Range("I3").Value = 0.1644281803937369
$ws.Range("J3").Value = 0.1663774220310514
$ws.Range("M3").Value = 3.248047
$ws.Range("N3").Value = 9.744140999999999
$ws.Range("O3").Value = 0.3236662202535691
$ws.Range("P3").Value = 0.3363578522370683
$ws.Range("Q3").Value = 15.68870254452967
$ws.Range("R3").Value = 141.198322900767
$ws.Range("S3").Value = 0.05321984765121283
$ws.Range("T3").Value = 0.05596235233510476
$ws.Range("G4").Value = 4.830195666666667
$ws.Range("H4").Value = 14.490587
$ws.Range("I4").Value = 0.1644281803937369
$ws.Range("J4").Value = 0.1663774220310514
$ws.Range("M4").Value = 1.1359575
$ws.Range("N4").Value = 2.271915
$ws.Range("O4").Value = 0.1131975831611099
$ws.Range("P4").Value = 0.07842419869182714
$ws.Range("Q4").Value = 5.4868969940175
$ws.Range("R4").Value = 32.921381964105
$ws.Range("S4").Value = 0.01861287262415002
$ws.Range("T4").Value = 0.01304801600319716
$ws.Range("G5").Value = 4.830195666666667
$ws.Range("H5").Value = 14.490587
$ws.Range("I5").Value = 0.1644281803937369
$ws.Range("J5").Value = 0.1663774220310514
$ws.Range("M5").Value = 5.637329
$ws.Range("N5").Value = 16.911987
$ws.Range("O5").Value = 0.5617569480231759
$ws.Range("P5").Value = 0.5837846172773179
$ws.Range("Q5").Value = 27.22940210737433
$ws.Range("R5").Value = 245.064618966369
$ws.Range("S5").Value = 0.09236867278698983
$ws.Range("T5").Value = 0.09712857964398416
$ws.Range("G6").Value = 4.674440333333334
$ws.Range("I6").Value = 0.1591260005621083
$ws.Range("J6").Value = 0.161012386613041
$ws.Range("M6").Value = 0.013841
$ws.Range("N6").Value = 0.041523
$ws.Range("O6").Value = 0.001379248562145083
$ws.Range("P6").Value = 0.001433331793786624
$ws.Range("Q6").Value = 0.06469892865366667
$ws.Range("R6").Value = 0.582290357883
$ws.Range("S6").Value = 0.0002194743074751854
$ws.Range("T6").Value = 0.0002307841729259355
$ws.Range("G7").Value = 4.674440333333334
$ws.Range("I7").Value = 0.1591260005621083
$ws.Range("J7").Value = 0.161012386613041
$ws.Range("M7").Value = 3.248047
$ws.Range("N7").Value = 9.744140999999999
$ws.Range("O7").Value = 0.3236662202535691
$ws.Range("P7").Value = 0.3363578522370683
$ws.Range("Q7").Value = 15.18280190136233
$ws.Range("R7").Value = 136.645217112261
$ws.Range("S7").Value = 0.05150371114600489
$ws.Range("T7").Value = 0.05415778054472697
$ws.Range("G8").Value = 4.674440333333334
$ws.Range("I8").Value = 0.1591260005621083
$ws.Range("J8").Value = 0.161012386613041
$ws.Range("M8").Value = 1.1359575
$ws.Range("N8").Value = 2.271915
$ws.Range("O8").Value = 0.1131975831611099
$ws.Range("P8").Value = 0.07842419869182714
$ws.Range("Q8").Value = 5.309965554952501
$ws.Range("R8").Value = 31.859793329715
$ws.Range("S8").Value = 0.01801267868172408
$ws.Range("T8").Value = 0.01262726739958642
$ws.Range("G9").Value = 4.674440333333334
$ws.Range("I9").Value = 0.1591260005621083
$ws.Range("J9").Value = 0.161012386613041
$ws.Range("M9").Value = 5.637329
$ws.Range("N9").Value = 16.911987
$ws.Range("O9").Value = 0.5617569480231759
$ws.Range("P9").Value = 0.5837846172773179
$ws.Range("Q9").Value = 26.35135804986967
$ws.Range("R9").Value = 237.162222448827
$ws.Range("S9").Value = 0.0893901364269041
$ws.Range("T9").Value = 0.09399655449580172
$ws.Range("G10").Value = 9.170097
$ws.Range("H10").Value = 27.510291
$ws.Range("I10").Value = 0.3121658971601493
$ws.Range("J10").Value = 0.3158665205145959
$ws.Range("M10").Value = 0.013841
$ws.Range("N10").Value = 0.041523
$ws.Range("O10").Value = 0.001379248562145083
$ws.Range("P10").Value = 0.001433331793786624
$ws.Range("Q10").Value = 0.126923312577
$ws.Range("R10").Value = 1.142309813193
$ws.Range("S10").Value = 0.0004305543648088657
$ws.Range("T10").Value = 0.0004527415264463252
$ws.Range("G11").Value = 9.170097
$ws.Range("H11").Value = 27.510291
$ws.Range("I11").Value = 0.3121658971601493
$ws.Range("J11").Value = 0.3158665205145959
$ws.Range("M11").Value = 3.248047
$ws.Range("N11").Value = 9.744140999999999
$ws.Range("O11").Value = 0.3236662202535691
$ws.Range("P11").Value = 0.3363578522370683
$ws.Range("Q11").Value = 29.784906050559
$ws.Range("R11").Value = 268.064154455031
$ws.Range("S11").Value = 0.1010375560258899
$ws.Range("T11").Value = 0.1062441844338854
$ws.Range("G12").Value = 9.170097
$ws.Range("H12").Value = 27.510291
$ws.Range("I12").Value = 0.3121658971601493
$ws.Range("J12").Value = 0.3158665205145959
$ws.Range("M12").Value = 1.1359575
$ws.Range("N12").Value = 2.271915
$ws.Range("O12").Value = 0.1131975831611099
$ws.Range("P12").Value = 0.07842419869182714
$ws.Range("Q12").Value = 10.4168404628775
$ws.Range("R12").Value = 62.50104277726501
$ws.Range("S12").Value = 0.0353364251038485
$ws.Range("T12").Value = 0.02477157876493276
$ws.Range("G13").Value = 9.170097
$ws.Range("H13").Value = 27.510291
$ws.Range("I13").Value = 0.3121658971601493
$ws.Range("J13").Value = 0.3158665205145959
$ws.Range("M13").Value = 5.637329
$ws.Range("N13").Value = 16.911987
$ws.Range("O13").Value = 0.5617569480231759
$ws.Range("P13").Value = 0.5837846172773179
$ws.Range("Q13").Value = 51.694853750913
$ws.Range("R13").Value = 465.2536837582171
$ws.Range("S13").Value = 0.1753613616656021
$ws.Range("T13").Value = 0.1843980157893315
$ws.Range("G14").Value = 1.0324785
$ws.Range("H14").Value = 2.064957
$ws.Range("I14").Value = 0.03514734656035429
$ws.Range("J14").Value = 0.02370933781115796
$ws.Range("M14").Value = 0.013841
$ws.Range("N14").Value = 0.041523
$ws.Range("O14").Value = 0.001379248562145083
$ws.Range("P14").Value = 0.001433331793786624
$ws.Range("Q14").Value = 0.0142905349185
$ws.Range("R14").Value = 0.08574320951100001
$ws.Range("S14").Value = 0.00004847692720658357
$ws.Range("T14").Value = 0.00003398334769436006
$ws.Range("G15").Value = 1.0324785
$ws.Range("H15").Value = 2.064957
$ws.Range("I15").Value = 0.03514734656035429
$ws.Range("J15").Value = 0.02370933781115796
$ws.Range("M15").Value = 3.248047
$ws.Range("N15").Value = 9.744140999999999
$ws.Range("O15").Value = 0.3236662202535691
$ws.Range("P15").Value = 0.3363578522370683
$ws.Range("Q15").Value = 3.3535386944895
$ws.Range("R15").Value = 20.121232166937
$ws.Range("S15").Value = 0.01137600881313216
$ws.Range("T15").Value = 0.007974821944124206
$ws.Range("G16").Value = 1.0324785
$ws.Range("H16").Value = 2.064957
$ws.Range("I16").Value = 0.03514734656035429
$ws.Range("J16").Value = 0.02370933781115796
$ws.Range("M16").Value = 1.1359575
$ws.Range("N16").Value = 2.271915
$ws.Range("O16").Value = 0.1131975831611099
$ws.Range("P16").Value = 0.07842419869182714
$ws.Range("Q16").Value = 1.17285169566375
$ws.Range("R16").Value = 4.691406782655
$ws.Range("S16").Value = 0.003978594685158057
$ws.Range("T16").Value = 0.001859385819353902
$ws.Range("G17").Value = 1.0324785
$ws.Range("H17").Value = 2.064957
$ws.Range("I17").Value = 0.03514734656035429
$ws.Range("J17").Value = 0.02370933781115796
$ws.Range("M17").Value = 5.637329
$ws.Range("N17").Value = 16.911987
$ws.Range("O17").Value = 0.5617569480231759
$ws.Range("P17").Value = 0.5837846172773179
$ws.Range("Q17").Value = 5.8204209899265
$ws.Range("R17").Value = 34.922525939559
$ws.Range("S17").Value = 0.0197442661348575
$ws.Range("T17").Value = 0.01384114669998549
$ws.Range("G18").Value = 9.668505333333334
$ws.Range("H18").Value = 29.005516
$ws.Range("I18").Value = 0.3291325753236513
$ws.Range("J18").Value = 0.3330343330301537
$ws.Range("M18").Value = 0.013841
$ws.Range("N18").Value = 0.041523
$ws.Range("O18").Value = 0.001379248562145083
$ws.Range("P18").Value = 0.001433331793786624
$ws.Range("Q18").Value = 0.1338217823186667
$ws.Range("R18").Value = 1.204396040868
$ws.Range("S18").Value = 0.0004539556312702541
$ws.Range("T18").Value = 0.0004773486979546421
$ws.Range("G19").Value = 9.668505333333334
$ws.Range("H19").Value = 29.005516
$ws.Range("I19").Value = 0.3291325753236513
$ws.Range("J19").Value = 0.3330343330301537
$ws.Range("M19").Value = 3.248047
$ws.Range("N19").Value = 9.744140999999999
$ws.Range("O19").Value = 0.3236662202535691
$ws.Range("P19").Value = 0.3363578522370683
$ws.Range("Q19").Value = 31.40375974241733
$ws.Range("R19").Value = 282.6338376817559
$ws.Range("S19").Value = 0.1065290966173293
$ws.Range("T19").Value = 0.112018712979227
$ws.Range("G20").Value = 9.668505333333334
$ws.Range("H20").Value = 29.005516
$ws.Range("I20").Value = 0.3291325753236513
$ws.Range("J20").Value = 0.3330343330301537
$ws.Range("M20").Value = 1.1359575
$ws.Range("N20").Value = 2.271915
$ws.Range("O20").Value = 0.1131975831611099
$ws.Range("P20").Value = 0.07842419869182714
$ws.Range("Q20").Value = 10.98301114719
$ws.Range("R20").Value = 65.89806688314
$ws.Range("S20").Value = 0.0372570120662293
$ws.Range("T20").Value = 0.0261179507047569
$ws.Range("G21").Value = 9.668505333333334
$ws.Range("H21").Value = 29.005516
$ws.Range("I21").Value = 0.3291325753236513
$ws.Range("J21").Value = 0.3330343330301537
$ws.Range("M21").Value = 5.637329
$ws.Range("N21").Value = 16.911987
$ws.Range("O21").Value = 0.5617569480231759
$ws.Range("P21").Value = 0.5837846172773179
$ws.Range("Q21").Value = 54.50454550225467
$ws.Range("R21").Value = 490.540909520292
$ws.Range("S21").Value = 0.1848925110088224
$ws.Range("T21").Value = 0.1944203206482151
